# This script normalizes the "Recorded By" list in column G so that the
# entry "System" (exact case) is always the first item in the comma-
# separated list, while the relative order of all other items is kept.
# Concretely: for any cell whose first item is not exactly "System" but
# which does contain "System" as one of its comma-separated items, the
# list is rotated left by one position (the first item is moved to the
# end of the list). Cells that already start with "System", that do not
# contain "System" at all, or that contain only a single item, are left
# untouched.

function Test-ExactEqual($s1, $s2) {
    # Ordinal / case-sensitive string comparison.
    # (The `-ceq`/`-cne` operators behave case-insensitively in this
    # runtime, so comparison is done manually via character codes.)
    if ($s1.Length -ne $s2.Length) { return $false }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        if ([int][char]$s1[$i] -ne [int][char]$s2[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Locate the "Recorded By" column from the header row, defaulting to
# column 7 (G) which is where it lives in this report layout.
$recordedByCol = 7
$headerRow = $ws.Cells.Item($firstRow, 1).Row
for ($c = 1; $c -le $used.Columns.Count; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Text
    if ($header -eq 'Recorded By') {
        $recordedByCol = $c
        break
    }
}

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ', '
    if ($parts.Count -le 1) { continue }

    $hasSystemExact = $false
    foreach ($p in $parts) {
        if (Test-ExactEqual $p 'System') { $hasSystemExact = $true }
    }
    $firstIsSystem = Test-ExactEqual $parts[0] 'System'

    if ($hasSystemExact -and -not $firstIsSystem) {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $newText = [string]::Join(', ', $rotated)
        $cell.Value = $newText
    }
}
